{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block (and the\n// blank paragraph immediately preceding it) that follows the last\n// bibliography entry ending in \"...Elsevier Science, New York, 1984.\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n\n// Find the index of the first target paragraph so we can also remove the\n// single blank paragraph that sits right before it (the spacer that was\n// separating the bibliography from the footer).\nlet firstTargetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (targets.indexOf(items[i].text.trim()) !== -1) {\n    firstTargetIndex = i;\n    break;\n  }\n}\n\nif (firstTargetIndex > 0 && items[firstTargetIndex - 1].text.trim() === \"\") {\n  items[firstTargetIndex - 1].delete();\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (targets.indexOf(text) !== -1) {\n    items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright footer block (and the\n# blank paragraph immediately preceding it) that follows the last\n# bibliography entry ending in \"...Elsevier Science, New York, 1984.\").\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Collect the indices (1-based, Word COM style) of the paragraphs to remove:\n# the two footer paragraphs plus the single blank spacer paragraph that\n# immediately precedes the first of them.\n$count = $d.Paragraphs.Count\n$toDelete = New-Object System.Collections.ArrayList\n\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($targets -contains $text) {\n        if ($toDelete.Count -eq 0 -and $i -gt 1) {\n            $prevText = $d.Paragraphs.Item($i - 1).Range.Text.Trim()\n            if ($prevText -eq \"\") {\n                [void]$toDelete.Add($i - 1)\n            }\n        }\n        [void]$toDelete.Add($i)\n    }\n}\n\n# Delete from the highest index to the lowest so earlier indices stay valid.\n$sorted = $toDelete | Sort-Object -Descending\nforeach ($idx in $sorted) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
